$wb = $excel.ActiveWorkbook

# Map of row -> new "想去人数" (F column) value that changed in this update
$updates = @{
    2  = 11802
    3  = 11514
    11 = 10839
    12 = 4188
    19 = 4
    23 = 10951
}

# Both "展览" and "全部类型" sheets contain identical data and both need the update
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
